# Update automatico via Actualizar 05-11-2020 17-45-31
#
# The source data refresh re-typed the "CodigoNiv1" column (F) as a
# zero-padded text code (matching the already-text "Admin1_id" style used
# elsewhere in the sheet) instead of a bare number, and backfilled the
# "Adm1tipo" column (L, "Municipio") for the newly added hospital rows
# further down the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    # Force a numeric-looking string ("08", "07", ...) to be stored as text
    # instead of being auto-coerced back into a number by Excel, while
    # restoring the cell's original number format afterwards so the visible
    # style/format of the cell does not otherwise change.
    $fmt = $cell.NumberFormat
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = $fmt
}

# --- Column F (CodigoNiv1): re-typed from plain numbers to zero-padded text codes ---
Set-TextValue $ws.Range("F84") "08"
Set-TextValue $ws.Range("F85") "08"
Set-TextValue $ws.Range("F87") "11"
Set-TextValue $ws.Range("F88") "05"
Set-TextValue $ws.Range("F89") "05"
Set-TextValue $ws.Range("F90") "05"
Set-TextValue $ws.Range("F91") "06"
Set-TextValue $ws.Range("F92") "01"

# --- Rows 86 & 96: same F re-typing, plus a top border (matching rows
#     120/121 and 132-135 which already use this "El Paraiso"/"Valle"
#     municipio code) picked up by copying the existing format. The text
#     re-type happens *before* the format paste so the border is the last
#     thing applied to the cell (avoids Excel minting a redundant style
#     that bundles the border with an explicitly-spelled-out "General"
#     number format). ---
Set-TextValue $ws.Range("F86") "07"
$ws.Range("E120:F120").Copy() | Out-Null
$ws.Range("E86:F86").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

Set-TextValue $ws.Range("F96") "17"
$ws.Range("E132:F132").Copy() | Out-Null
$ws.Range("E96:F96").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# --- Column L (Adm1tipo): backfill "Municipio" for rows 97-135 ---
for ($r = 97; $r -le 135; $r++) {
    $ws.Range("L$r").Value = "Municipio"
}

# --- View/formatting touch-ups from the refresh ---
$ws.Columns.Item(15).ColumnWidth = 31.140625
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("L2").Select() | Out-Null
$ws.Range("E86:F86").Select() | Out-Null
